$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.892.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.236.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.04%  "

$ws.Range("E4").Value = "  +0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.19%  "

$ws.Range("E7").Value = "  -3.19%  "

$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("E9").Value = "  -7.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0820"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.90%  "

$ws.Range("E13").Value = "  -2.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.578.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.240.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.838"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.66%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.752.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0971"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.48%  "

$ws.Range("E24").Value = "  -7.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.28%  "

$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0829"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.83%  "

$ws.Range("E35").Value = "  -3.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.40%  "

$ws.Range("E37").Value = "  -6.48%  "

$ws.Range("E38").Value = "  -3.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -13.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -11.40%  "

$ws.Range("E42").Value = "  -7.32%  "

$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.705.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "82.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.51%  "

$ws.Range("E46").Value = "  -7.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "71.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "56.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.29%  "

